# issue #5: add legislator_id, name, date into dataframe
#
# The 股票 (stocks) sheet gains three new trailing columns - date,
# legislator_name, legislator_id - with the same value repeated for
# every existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count          # includes header row 1
$lastCol = $usedRange.Columns.Count       # currently G (7)

$dateCol = $lastCol + 1                   # H
$nameCol = $lastCol + 2                   # I
$idCol   = $lastCol + 3                   # J

# --- Header row: style like the existing header cells (bold + border). ---
$ws.Cells.Item(1, $lastCol).Copy()
$ws.Range($ws.Cells.Item(1, $dateCol), $ws.Cells.Item(1, $idCol)).PasteSpecial(-4122)

$ws.Cells.Item(1, $dateCol).Value = "date"
$ws.Cells.Item(1, $nameCol).Value = "legislator_name"
$ws.Cells.Item(1, $idCol).Value = "legislator_id"

# --- Data rows: date (kept as text, not a date serial), legislator name,
# and legislator id for each existing row. ---
$ws.Range($ws.Cells.Item(2, $dateCol), $ws.Cells.Item($lastRow, $dateCol)).NumberFormat = "@"

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $dateCol).Value = "2011-11-21"
    $ws.Cells.Item($r, $nameCol).Value = "蔣乃辛"
    $ws.Cells.Item($r, $idCol).Value = 1722
}

# Drop the text-number-format marker again so the new data cells render
# with the same plain look as the rest of the sheet.
$ws.Range($ws.Cells.Item(2, $dateCol), $ws.Cells.Item($lastRow, $dateCol)).Style = "Normal"
